$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reference field corrections (BOM reference designators) ---
$ws.Range("B4").Value = "C7, C9, C11, C12"
$ws.Range("B5").Value = "C8, C10"
$ws.Range("B10").Value = "R1-R62"

# --- Update U1 pricing / note, and restyle the price cell ---
# Copy the number/font formatting used elsewhere in the Mouser P/N column (fontId 3)
# onto I12, then override number format (currency, built-in id 8) and remove the border
# so the resulting style matches the new target style exactly.
$ws.Range("F4").Copy()
$ws.Range("I12").PasteSpecial(-4122)
$ws.Range("I12").NumberFormat = '"$"#,##0.00;[Red]\-"$"#,##0.00'
$ws.Range("I12").Borders.LineStyle = -4142
$ws.Range("I12").Value = 40.27
$ws.Range("K12").Value = "$42.39 From Digikey"

# --- Remove the now-obsolete "U1 Alt" row (row 13) contents, keep formatting/styles ---
$ws.Range("B13:K13").ClearContents()

# --- Cosmetic column width adjustments (closest achievable values) ---
$ws.Columns("B").ColumnWidth = 13.85
$ws.Columns("I").ColumnWidth = 8.6

# --- Restore the active selection cell as recorded in the saved view ---
$ws.Range("K17").Select()
